$wb = $excel.ActiveWorkbook

$wsIntegration = $wb.Worksheets.Item("Integration")
$wsSystem = $wb.Worksheets.Item("System")

# ---------------------------------------------------------------------------
# 1. Add the new "Functional" worksheet after "System".
# ---------------------------------------------------------------------------
$wsFunctional = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsSystem)
$wsFunctional.Name = "Functional"

# Copy the header row formatting (fill + border + centered alignment) from the
# Integration sheet's header row so the new sheet matches the existing look.
$wsIntegration.Range("A1:G1").Copy()
$wsFunctional.Range("A1:Y1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Fill in the header row and the data row (row 2) for the new product test
#    case. The write order below matches the original authoring order so
#    that newly interned shared strings land at the same table positions.
# ---------------------------------------------------------------------------
$wsFunctional.Range("A1").Value = "TC_ID"
$wsFunctional.Range("B1").Value = "TC_Name"

$wsFunctional.Range("A2").Value = "TC01_Products"
$wsFunctional.Range("B2").Value = "AddProduct"

$wsFunctional.Range("C1").Value = "ProductName"
$wsFunctional.Range("C2").Value = "Product_123"

$wsFunctional.Range("D1").Value = "Quantity"
$wsFunctional.Range("D2").NumberFormat = "@"
$wsFunctional.Range("D2").Value = "10"

$wsFunctional.Range("E1").Value = "Price"
$wsFunctional.Range("F1").Value = "Category"

$wsFunctional.Range("E2").NumberFormat = "@"
$wsFunctional.Range("E2").Value = "15000"

$wsFunctional.Range("F2").Value = "Electronics"

$wsFunctional.Range("G1").Value = "Select Vendor"

# Column widths to roughly match the authored layout.
$wsFunctional.Columns.Item(1).ColumnWidth = 16
$wsFunctional.Columns.Item(2).ColumnWidth = 22
$wsFunctional.Columns.Item(3).ColumnWidth = 13.5
$wsFunctional.Columns.Item(4).ColumnWidth = 16
$wsFunctional.Columns.Item(5).ColumnWidth = 13.9
$wsFunctional.Columns.Item(6).ColumnWidth = 17.6
$wsFunctional.Columns.Item(7).ColumnWidth = 23.5
$wsFunctional.Columns.Item(8).ColumnWidth = 11.3
$wsFunctional.Columns.Item(10).ColumnWidth = 12.6

# ---------------------------------------------------------------------------
# 4. Update the "System" sheet's Purchase-Order test row (row 2) so it
#    references the new product/contact records.
# ---------------------------------------------------------------------------
$wsSystem.Range("K2").NumberFormat = "@"
$wsSystem.Range("K2").Value = "75000"

$wsSystem.Range("Y2").NumberFormat = "@"
$wsSystem.Range("Y2").Value = "568989"

$wsSystem.Range("S2").NumberFormat = "@"
$wsSystem.Range("S2").Value = "564443"

$wsSystem.Range("E2").NumberFormat = "@"
$wsSystem.Range("E2").Value = "9874234789"

$wsSystem.Range("L2").Value = "VID_002"

$wsSystem.Range("N2").NumberFormat = "@"
$wsSystem.Range("N2").Value = "25-06-2025"

$wsSystem.Range("J2").NumberFormat = "@"
$wsSystem.Range("J2").Value = "10"

$wsSystem.Range("AA1").Value = "TargetSize"

$wsSystem.Range("AA2").NumberFormat = "@"
$wsSystem.Range("AA2").Value = "7"

# ---------------------------------------------------------------------------
# 5. Back on the Functional sheet: set the vendor id, carrying over the grey
#    font color that came along with the pasted vendor id text.
# ---------------------------------------------------------------------------
$wsFunctional.Range("G2").Font.Color = 2039583
$wsFunctional.Range("G2").Value = "VID_002"

# ---------------------------------------------------------------------------
# 6. Page setup - portrait orientation for the updated/added sheets.
# ---------------------------------------------------------------------------
$wsSystem.PageSetup.Orientation = 1
$wsFunctional.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 7. Restore/update the on-screen selection state for each sheet.
# ---------------------------------------------------------------------------
$wsIntegration.Activate()
$wsIntegration.Range("A1:XFD1").Select()

$wsFunctional.Activate()
$wsFunctional.Range("G2").Select()

$wsSystem.Activate()
$wsSystem.Range("N6").Select()
